# Rename the sheet "casa_mia" -> "home_buccinasco".
# Excel automatically re-points every formula / defined name (solver_*,
# etc.) that referenced the old sheet name, exactly like the OOXML diff
# shows for xl/workbook.xml.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("casa_mia")
$ws.Name = "home_buccinasco"

# --- Row 80: the capacity-units label loses its border formatting, and
# the (empty, but styled) cells to its right are deleted outright. ---
$ws.Range("AL80:AX80").Clear()
$ws.Range("AK80").ClearFormats()

# --- Row 81: the data row keeps its fill but loses its border. ---
$ws.Range("AK81:AX81").Borders.LineStyle = -4142   # xlLineStyleNone

# --- View state: scroll position / active selection changed. ---
$ws.Activate()
$ws.Range("O14").Select()
